# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.813.49"
$ws.Range("E2").Value = "'  -0.46%  "
$ws.Range("D3").Value = "'1.635.64"
$ws.Range("E3").Value = "'  -0.10%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'215.67"
$ws.Range("E5").Value = "'  +0.65%  "
$ws.Range("D6").Value = "'0.5054"
$ws.Range("E6").Value = "'  -0.22%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'0.2572"
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.06416"
$ws.Range("E9").Value = "'  +0.93%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("E10").Value = "'  -1.05%  "
$ws.Range("D11").Value = "'0.07770"
$ws.Range("E11").Value = "'  +0.42%  "
$ws.Range("D12").Value = "'4.274"
$ws.Range("E12").Value = "'  -0.45%  "
$ws.Range("D13").Value = "'1.862.96"
$ws.Range("E13").Value = "'  +0.04%  "
$ws.Range("D14").Value = "'1.631.23"
$ws.Range("E14").Value = "'  -1.00%  "
$ws.Range("D15").Value = "'0.5627"
$ws.Range("E15").Value = "'  +3.16%  "
$ws.Range("B16").Value = "'Litecoin"
$ws.Range("C16").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'63.12"
$ws.Range("E16").Value = "'  -1.58%  "
$ws.Range("B17").Value = "'ShibaInu"
$ws.Range("C17").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0₅7587"
$ws.Range("E17").Value = "'  -1.90%  "
$ws.Range("D18").Value = "'25.847.40"
$ws.Range("E18").Value = "'  -0.37%  "
$ws.Range("D20").Value = "'195.19"
$ws.Range("E20").Value = "'  -0.35%  "
$ws.Range("D21").Value = "'4.320"
$ws.Range("E21").Value = "'  -2.97%  "
$ws.Range("D22").Value = "'9.871"
$ws.Range("E22").Value = "'  -0.62%  "
$ws.Range("D23").Value = "'6.088"
$ws.Range("E23").Value = "'  -0.87%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "'  -0.02%  "
$ws.Range("E25").Value = "'  -5.26%  "
$ws.Range("D26").Value = "'0.1271"
$ws.Range("E26").Value = "'  +1.90%  "
$ws.Range("D27").Value = "'139.70"
$ws.Range("E27").Value = "'  -2.30%  "
$ws.Range("D28").Value = "'6.779"
$ws.Range("E28").Value = "'  -0.93%  "
$ws.Range("D29").Value = "'15.47"
$ws.Range("E29").Value = "'  -0.70%  "
$ws.Range("D30").Value = "'1.241"
$ws.Range("E30").Value = "'  +0.38%  "
$ws.Range("D31").Value = "'0.04873"
$ws.Range("E31").Value = "'  -0.12%  "
$ws.Range("D32").Value = "'3.292"
$ws.Range("E32").Value = "'  +1.40%  "
$ws.Range("D33").Value = "'3.213"
$ws.Range("E33").Value = "'  +0.34%  "
$ws.Range("D34").Value = "'1.554"
$ws.Range("E34").Value = "'  +0.15%  "
$ws.Range("D35").Value = "'2.367"
$ws.Range("E35").Value = "'  -0.10%  "
$ws.Range("D36").Value = "'0.9021"
$ws.Range("E36").Value = "'  -1.20%  "
$ws.Range("D37").Value = "'2.577"
$ws.Range("E37").Value = "'  +0.25%  "
$ws.Range("D38").Value = "'1.130.18"
$ws.Range("E38").Value = "'  +0.53%  "
$ws.Range("D39").Value = "'0.5499"
$ws.Range("E39").Value = "'  -0.22%  "
$ws.Range("D40").Value = "'0.01562"
$ws.Range("E40").Value = "'  -0.24%  "
$ws.Range("D41").Value = "'0.9951"
$ws.Range("E41").Value = "'  -0.61%  "
$ws.Range("D42").Value = "'5.516"
$ws.Range("E42").Value = "'  -1.57%  "
$ws.Range("D43").Value = "'0.8004"
$ws.Range("E43").Value = "'  -0.43%  "
$ws.Range("D44").Value = "'97.81"
$ws.Range("E44").Value = "'  -0.70%  "
$ws.Range("D45").Value = "'1.773.39"
$ws.Range("E45").Value = "'  +0.07%  "
$ws.Range("D46").Value = "'0.0₈114"
$ws.Range("E46").Value = "'  -6.25%  "
$ws.Range("D47").Value = "'55.31"
$ws.Range("D48").Value = "'0.4403"
$ws.Range("E48").Value = "'  -1.94%  "
$ws.Range("D49").Value = "'0.05050"
$ws.Range("E49").Value = "'  -2.52%  "
$ws.Range("D50").Value = "'7.669"
$ws.Range("E50").Value = "'  +1.83%  "
$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "'  +0.10%  "

# Reset style on touched cells so no stray quotePrefix style is left applied
$touched = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "D18", "E18", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $touched) { $ws.Range($addr).Style = "Normal" }
